$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header label
$ws.Range("E1").Value = "Memory Usage (mb)"

# Update Run Time (ms) column D
$ws.Range("D2").Value = 22.47881889343262
$ws.Range("D3").Value = 18.53322982788086
$ws.Range("D4").Value = 18.04685592651367
$ws.Range("D5").Value = 18.06902885437012
$ws.Range("D6").Value = 18.32818984985352

# Update Memory Usage column E
$ws.Range("E2").Value = 62.324736
$ws.Range("E3").Value = 62.357504
$ws.Range("E4").Value = 62.357504
$ws.Range("E5").Value = 62.357504
$ws.Range("E6").Value = 62.357504
